$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 174.5
$ws.Range("I29").Value = 174.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 523.5
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -242.5

$ws.Range("H43").Value = 1694.0769
$ws.Range("J43").Value = 1407.6364
$ws.Range("L43").Value = 1407.6364
$ws.Range("N43").Value = -1545.6364

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws.Range("H86").Value = 4269.8096
$ws.Range("I86").Value = 1461.25
$ws.Range("J86").Value = 5998.154
$ws.Range("K86").Value = 1461.25
$ws.Range("L86").Value = 5998.154
$ws.Range("M86").Value = -338.25
$ws.Range("N86").Value = -8244.154

$ws.Range("H89").Value = 4269.8096
$ws.Range("I89").Value = 1461.25
$ws.Range("J89").Value = 5998.154
$ws.Range("K89").Value = 7306.25
$ws.Range("L89").Value = 29990.77
$ws.Range("M89").Value = -1690.25
$ws.Range("N89").Value = -41222.77

$ws.Range("H126").Value = 48000
$ws.Range("J126").Value = 48000
$ws.Range("L126").Value = 48000
$ws.Range("N126").Value = -57880

$ws.Range("H129").Value = 438124.5
$ws.Range("I129").Value = 17049.5
$ws.Range("J129").Value = 535295.6
$ws.Range("K129").Value = 51148.5
$ws.Range("L129").Value = 1605886.8
$ws.Range("M129").Value = -46148.5
$ws.Range("N129").Value = -1615886.8

$ws.Range("H130").Value = 50390
$ws.Range("J130").Value = 50390
$ws.Range("L130").Value = 50390
$ws.Range("N130").Value = -60430

$ws.Range("H138").Value = 2598.962
$ws.Range("I138").Value = 1369.7858
$ws.Range("J138").Value = 3273.804
$ws.Range("K138").Value = 4109.357400000001
$ws.Range("L138").Value = 9821.412
$ws.Range("M138").Value = 1030.642599999999
$ws.Range("N138").Value = -20101.412

$ws.Range("H140").Value = 50251.668
$ws.Range("J140").Value = 50251.668
$ws.Range("L140").Value = 50251.668
$ws.Range("N140").Value = -60611.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 252495
$ws.Range("I2").Value = 3240
$ws.Range("K2").Value = 3240
$ws.Range("M2").Value = -3127

$ws.Range("H32").Value = 19360.18
$ws.Range("I32").Value = 3296.988
$ws.Range("K32").Value = 3296.988
$ws.Range("M32").Value = -3009.988

$ws.Range("H61").Value = 2224.842
$ws.Range("I61").Value = 1722.6316
$ws.Range("J61").Value = 2727.0527
$ws.Range("K61").Value = 1722.6316
$ws.Range("L61").Value = 2727.0527
$ws.Range("M61").Value = -1510.6316
$ws.Range("N61").Value = -3151.0527

$ws.Range("H88").Value = 2908
$ws.Range("I88").Value = 3100
$ws.Range("J88").Value = 2825.7144
$ws.Range("K88").Value = 3100
$ws.Range("L88").Value = 2825.7144
$ws.Range("M88").Value = -2694
$ws.Range("N88").Value = -3637.7144

$ws.Range("H91").Value = 2908
$ws.Range("I91").Value = 3100
$ws.Range("J91").Value = 2825.7144
$ws.Range("K91").Value = 3100
$ws.Range("L91").Value = 2825.7144
$ws.Range("M91").Value = -1696
$ws.Range("N91").Value = -5633.7144

$ws.Range("H116").Value = 252495
$ws.Range("I116").Value = 3240
$ws.Range("K116").Value = 3240
$ws.Range("M116").Value = -946

$ws.Range("H136").Value = 2224.842
$ws.Range("I136").Value = 1722.6316
$ws.Range("J136").Value = 2727.0527
$ws.Range("K136").Value = 5167.8948
$ws.Range("L136").Value = 8181.158100000001
$ws.Range("M136").Value = -2617.8948
$ws.Range("N136").Value = -13281.1581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 252495
$ws.Range("I3").Value = 3240
$ws.Range("K3").Value = 3240
$ws.Range("M3").Value = -3126

$ws.Range("H20").Value = 57080.945
$ws.Range("I20").Value = 68163.2
$ws.Range("K20").Value = 68163.2
$ws.Range("M20").Value = -67916.2

$ws.Range("H86").Value = 52853.457
$ws.Range("I86").Value = 73064.47
$ws.Range("J86").Value = 3769.5715
$ws.Range("K86").Value = 73064.47
$ws.Range("L86").Value = 3769.5715
$ws.Range("M86").Value = -71941.47
$ws.Range("N86").Value = -6015.5715

$ws.Range("H89").Value = 52853.457
$ws.Range("I89").Value = 73064.47
$ws.Range("J89").Value = 3769.5715
$ws.Range("K89").Value = 365322.35
$ws.Range("L89").Value = 18847.8575
$ws.Range("M89").Value = -359706.35
$ws.Range("N89").Value = -30079.8575

$ws.Range("H107").Value = 83401160
$ws.Range("I107").Value = 90982904
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 90982904
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -90980984
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 24142.5
$ws.Range("I19").Value = 310
$ws.Range("J19").Value = 47975
$ws.Range("K19").Value = 310
$ws.Range("L19").Value = 47975
$ws.Range("M19").Value = -140
$ws.Range("N19").Value = -48315

$ws.Range("H24").Value = 24142.5
$ws.Range("I24").Value = 310
$ws.Range("J24").Value = 47975
$ws.Range("K24").Value = 310
$ws.Range("L24").Value = 47975
$ws.Range("M24").Value = -140
$ws.Range("N24").Value = -48315

$ws.Range("H31").Value = 38492.63
$ws.Range("I31").Value = 1294.7142
$ws.Range("J31").Value = 73999.73
$ws.Range("K31").Value = 1294.7142
$ws.Range("L31").Value = 73999.73
$ws.Range("M31").Value = -999.7141999999999
$ws.Range("N31").Value = -74589.73

$ws.Range("H34").Value = 38492.63
$ws.Range("I34").Value = 1294.7142
$ws.Range("J34").Value = 73999.73
$ws.Range("K34").Value = 1294.7142
$ws.Range("L34").Value = 73999.73
$ws.Range("M34").Value = -1092.7142
$ws.Range("N34").Value = -74403.73

$ws.Range("H92").Value = 23333
$ws.Range("J92").Value = 23333
$ws.Range("L92").Value = 23333
$ws.Range("N92").Value = -28325

$ws.Range("H132").Value = 3667.6
$ws.Range("I132").Value = 3580.3872
$ws.Range("K132").Value = 10741.1616
$ws.Range("M132").Value = -8211.161599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 95.666664
$ws.Range("J12").Value = 95.666664
$ws.Range("L12").Value = 286.999992
$ws.Range("N12").Value = -632.999992

$ws.Range("H33").Value = 1129.3684
$ws.Range("J33").Value = 1912.909
$ws.Range("L33").Value = 11477.454
$ws.Range("N33").Value = -12043.454

$ws.Range("H34").Value = 671.6667
$ws.Range("I34").Value = 123.333336
$ws.Range("J34").Value = 945.8333
$ws.Range("K34").Value = 370.000008
$ws.Range("L34").Value = 2837.4999
$ws.Range("M34").Value = -286.000008
$ws.Range("N34").Value = -3005.4999

$ws.Range("H131").Value = 803.79
$ws.Range("J131").Value = 848.5
$ws.Range("L131").Value = 2545.5
$ws.Range("N131").Value = -12625.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 77782
$ws.Range("I70").Value = 101132.86
$ws.Range("J70").Value = 7729.4287
$ws.Range("K70").Value = 101132.86
$ws.Range("L70").Value = 7729.4287
$ws.Range("M70").Value = -100862.86
$ws.Range("N70").Value = -8269.4287

$ws.Range("H73").Value = 77782
$ws.Range("I73").Value = 101132.86
$ws.Range("J73").Value = 7729.4287
$ws.Range("K73").Value = 101132.86
$ws.Range("L73").Value = 7729.4287
$ws.Range("M73").Value = -100196.86
$ws.Range("N73").Value = -9601.4287

$ws.Range("H102").Value = 1954.9166
$ws.Range("I102").Value = 1954.9166
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1954.9166
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -332.9166

$ws.Range("H122").Value = 1435.1666
$ws.Range("I122").Value = 1513.421
$ws.Range("J122").Value = 1137.8
$ws.Range("K122").Value = 4540.263
$ws.Range("L122").Value = 3413.4
$ws.Range("M122").Value = -2090.263
$ws.Range("N122").Value = -8313.4

$ws.Range("H127").Value = 48000
$ws.Range("J127").Value = 48000
$ws.Range("L127").Value = 48000
$ws.Range("N127").Value = -57920

$ws.Range("H132").Value = 3460.7896
$ws.Range("I132").Value = 3150.5386
$ws.Range("K132").Value = 9451.6158
$ws.Range("M132").Value = -6921.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2123.72
$ws.Range("I7").Value = 1445.9333
$ws.Range("J7").Value = 3140.4
$ws.Range("K7").Value = 1445.9333
$ws.Range("L7").Value = 3140.4
$ws.Range("M7").Value = -1333.9333
$ws.Range("N7").Value = -3364.4

$ws.Range("H40").Value = 75197.28999999999
$ws.Range("I40").Value = 172043.67
$ws.Range("K40").Value = 172043.67
$ws.Range("M40").Value = -171907.67

$ws.Range("H126").Value = 2123.72
$ws.Range("I126").Value = 1445.9333
$ws.Range("J126").Value = 3140.4
$ws.Range("K126").Value = 4337.7999
$ws.Range("L126").Value = 9421.200000000001
$ws.Range("M126").Value = -1867.7999
$ws.Range("N126").Value = -14361.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3739.3333
$ws.Range("I132").Value = 3421.889
$ws.Range("J132").Value = 4374.222
$ws.Range("K132").Value = 10265.667
$ws.Range("L132").Value = 13122.666
$ws.Range("M132").Value = -18182.666
